$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text would otherwise be auto-parsed as a number by Excel;
# force them to remain plain text, matching the original inlineStr cell type,
# then restore the default (unstyled) cell style so no stray formatting is left behind.
$textCells = @('D4', 'D5', 'D6', 'D7', 'D8', 'D9', 'D10', 'D11', 'D12', 'D14', 'D15', 'D16', 'D18', 'D21', 'D23', 'D24', 'D26', 'D27', 'D28', 'D29', 'D30', 'D31', 'D32', 'D33', 'D34', 'D35', 'D36', 'D37', 'D39', 'D40', 'D41', 'D42', 'D43', 'D45', 'D46', 'D47', 'D49', 'D50', 'D51')
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = '@'
}

$ws.Range('D2').Value = '29.866.17'
$ws.Range('E2').Value = '  +2.13%  '
$ws.Range('D3').Value = '1.872.07'
$ws.Range('E3').Value = '  +0.84%  '
$ws.Range('D4').Value = '0.9998'
$ws.Range('E4').Value = '  -0.03%  '
$ws.Range('D5').Value = '247.18'
$ws.Range('E5').Value = '  +2.38%  '
$ws.Range('D6').Value = '0.7039'
$ws.Range('E6').Value = '  +1.48%  '
$ws.Range('D7').Value = '1.000'
$ws.Range('D8').Value = '0.07800'
$ws.Range('E8').Value = '  +0.47%  '
$ws.Range('D9').Value = '0.3102'
$ws.Range('E9').Value = '  +1.32%  '
$ws.Range('D10').Value = '24.06'
$ws.Range('E10').Value = '  +1.32%  '
$ws.Range('D11').Value = '0.07859'
$ws.Range('E11').Value = '  +0.65%  '
$ws.Range('D12').Value = '5.203'
$ws.Range('E12').Value = '  +1.90%  '
$ws.Range('D13').Value = '1.870.69'
$ws.Range('E13').Value = '  +0.57%  '
$ws.Range('D14').Value = '93.44'
$ws.Range('E14').Value = '  +1.31%  '
$ws.Range('D15').Value = '0.6997'
$ws.Range('E15').Value = '  +1.83%  '
$ws.Range('D16').Value = '6.650'
$ws.Range('E16').Value = '  +1.93%  '
$ws.Range('D17').Value = '29.870.21'
$ws.Range('E17').Value = '  +2.11%  '
$ws.Range('D18').Value = '0.000008442'
$ws.Range('E18').Value = '  +0.57%  '
$ws.Range('D19').Value = '244.88'
$ws.Range('E19').Value = '  -1.02%  '
$ws.Range('D20').Value = '2.111.33'
$ws.Range('E20').Value = '  -0.07%  '
$ws.Range('D21').Value = '12.89'
$ws.Range('E21').Value = '  +0.51%  '
$ws.Range('E22').Value = '  -0.01%  '
$ws.Range('D23').Value = '7.679'
$ws.Range('E23').Value = '  +1.62%  '
$ws.Range('D24').Value = '1.000'
$ws.Range('E24').Value = '  -0.04%  '
$ws.Range('D25').Value = '0.1522'
$ws.Range('E25').Value = '  +1.84%  '
$ws.Range('D26').Value = '9.007'
$ws.Range('E26').Value = '  +1.63%  '
$ws.Range('D27').Value = '160.60'
$ws.Range('E27').Value = '  -0.65%  '
$ws.Range('D28').Value = '18.50'
$ws.Range('E28').Value = '  +0.27%  '
$ws.Range('D29').Value = '1.550'
$ws.Range('E29').Value = '  -0.08%  '
$ws.Range('D30').Value = '4.310'
$ws.Range('E30').Value = '  +1.20%  '
$ws.Range('D31').Value = '4.262'
$ws.Range('E31').Value = '  +1.65%  '
$ws.Range('D32').Value = '1.205'
$ws.Range('E32').Value = '  +0.61%  '
$ws.Range('D33').Value = '0.05149'
$ws.Range('E33').Value = '  -0.36%  '
$ws.Range('D34').Value = '0.7906'
$ws.Range('E34').Value = '  +3.29%  '
$ws.Range('D35').Value = '1.941'
$ws.Range('E35').Value = '  +5.40%  '
$ws.Range('D36').Value = '1.170'
$ws.Range('E36').Value = '  +0.18%  '
$ws.Range('D37').Value = '2.708'
$ws.Range('D38').Value = '1.338.78'
$ws.Range('E38').Value = '  +8.85%  '
$ws.Range('D39').Value = '0.01895'
$ws.Range('E39').Value = '  +2.22%  '
$ws.Range('D40').Value = '2.751'
$ws.Range('E40').Value = '  +0.80%  '
$ws.Range('D41').Value = '0.9566'
$ws.Range('E41').Value = '  +4.57%  '
$ws.Range('D42').Value = '6.097'
$ws.Range('E42').Value = '  +11.67%  '
$ws.Range('D43').Value = '108.01'
$ws.Range('E43').Value = '  -1.40%  '
$ws.Range('E44').Value = '  -0.02%  '
$ws.Range('B45').Value = 'BabyDogeCoin'
$ws.Range('C45').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D45').Value = '0.00000000126'
$ws.Range('E45').Value = '  +2.23%  '
$ws.Range('B46').Value = 'EnergySwap'
$ws.Range('C46').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D46').Value = '9.865'
$ws.Range('E46').Value = '  +3.67%  '
$ws.Range('B47').Value = 'Aave'
$ws.Range('C47').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D47').Value = '66.09'
$ws.Range('E47').Value = '  +2.39%  '
$ws.Range('B48').Value = 'RocketPoolETH'
$ws.Range('C48').Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range('D48').Value = '2.009.26'
$ws.Range('E48').Value = '  +0.14%  '
$ws.Range('D49').Value = '0.5213'
$ws.Range('E49').Value = '  +0.54%  '
$ws.Range('D50').Value = '1.793'
$ws.Range('E50').Value = '  +2.76%  '
$ws.Range('D51').Value = '7.059'
$ws.Range('E51').Value = '  +0.89%  '

foreach ($addr in $textCells) {
    $ws.Range($addr).Style = 'Normal'
}
